# tax_investments_helper.xlsx - "Add files via upload"
#
# 1) D48 had a leftover/incorrect label ("upto 4 lacs"); fix it to read
#    "3 - 5 lacs" so it lines up with the "20 - 30 lacs per annum" example
#    income row it sits next to.
# 2) Append a brand-new two-column summary table in rows 54-60:
#    "Component" / "Max Tax Investments (under old regime)" header, five
#    component rows, and a "Total" footer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the mislabeled example-income row ---------------------------
$ws.Range("D48").Value2 = "3 - 5 lacs"

# --- 2) New "Max Tax Investments" table ----------------------------------
$ws.Range("C54").Value2 = "Component"
$ws.Range("D54").Value2 = "Max Tax Investments (under old regime)"

$ws.Range("C55").Value2 = "PPFO"
$ws.Range("D55").Value2 = "1.5L"

$ws.Range("C56").Value2 = "NPS"
$ws.Range("D56").Value2 = "0.5L"

$ws.Range("C57").Value2 = "Rent"
$ws.Range("D57").Value2 = "6.0L - 12.0L"

$ws.Range("C58").Value2 = "MI"
$ws.Range("D58").Value2 = "0.5L"

$ws.Range("C59").Value2 = "HL"
$ws.Range("D59").Value2 = "2L"

$ws.Range("C60").Value2 = "Total"
$ws.Range("D60").Value2 = "10.5L - 16.5L"

# Match the look of the existing tables: bold/bordered header+total rows
# (style of row 44, the "Example Income" header) and plain bordered data
# rows (style of row 48) for the body rows in between.
$ws.Range("C44:D44").Copy()
$ws.Range("C54:D54").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C48:D48").Copy()
$ws.Range("C55:D59").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C44:D44").Copy()
$ws.Range("C60:D60").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Match the saved view/selection state --------------------------------
$ws.Activate()
$ws.Range("D29").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
